# Updates cryptos list values (price + 1h volume change) per the Oct 16 2024 refresh,
# including the Bittensor/WrappedeETH row-28/29 re-rank swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value.
# Values that look like plain decimals (e.g. "602.08") are forced to Text format
# first so Excel does not reinterpret them as floating-point numbers (which would
# both change the cell type and introduce binary floating-point rounding noise),
# matching the source data which stores these as text.
$updates = [ordered]@{
    "D2" = "67.670.44"
    "E2" = "  +1.64%  "
    "D3" = "2.618.33"
    "E3" = "  +1.07%  "
    "E4" = "  -0.08%  "
    "D5" = "602.08"
    "E5" = "  +1.94%  "
    "D6" = "154.48"
    "E6" = "  +0.67%  "
    "E7" = "  +0.04%  "
    "E8" = "  +1.70%  "
    "D9" = "2.615.51"
    "E9" = "  +0.99%  "
    "E10" = "  +11.07%  "
    "E11" = "  +0.85%  "
    "D12" = "5.25"
    "E12" = "  +1.23%  "
    "E13" = "  -0.15%  "
    "D14" = "27.99"
    "E14" = "  -0.18%  "
    "E15" = "  +4.28%  "
    "D16" = "3.094.65"
    "D17" = "67.567.57"
    "E17" = "  +1.74%  "
    "D18" = "2.622.82"
    "E18" = "  +1.32%  "
    "D19" = "11.28"
    "E19" = "  +0.26%  "
    "D20" = "363.50"
    "E20" = "  +3.26%  "
    "D21" = "7.62"
    "E21" = "  -2.58%  "
    "E22" = "  -0.25%  "
    "D23" = "2.14"
    "E23" = "  +5.92%  "
    "E24" = "  +0.03%  "
    "D25" = "70.08"
    "E25" = "  +3.66%  "
    "D26" = "10.11"
    "E26" = "  -1.86%  "
    "E27" = "  +4.26%  "
    "B28" = "Bittensor"
    "C28" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D28" = "592.99"
    "E28" = "  +0.92%  "
    "B29" = "WrappedeETH"
    "C29" = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
    "D29" = "2.744.64"
    "E29" = "  +1.37%  "
    "D30" = "1.02"
    "E30" = "  +1.93%  "
    "E31" = "  -0.02%  "
    "D32" = "7.95"
    "E32" = "  -0.41%  "
    "E33" = "  +0.79%  "
    "E34" = "  -2.13%  "
    "E35" = "  +0.04%  "
    "D36" = "1.54"
    "E36" = "  -0.49%  "
    "D37" = "4.98"
    "E37" = "  +0.25%  "
    "D38" = "19.44"
    "E38" = "  +1.32%  "
    "D39" = "156.50"
    "E39" = "  +2.76%  "
    "D40" = "0.372"
    "E40" = "  +1.10%  "
    "D41" = "5.43"
    "E41" = "  +0.53%  "
    "D42" = "1.86"
    "E42" = "  +3.76%  "
    "D43" = "2.66"
    "E43" = "  +4.02%  "
    "E44" = "  -0.28%  "
    "E45" = "  +0.08%  "
    "D46" = "16.43"
    "E46" = "  +0.05%  "
    "D47" = "156.99"
    "E47" = "  +0.60%  "
    "D48" = "0.0₆0291"
    "E48" = "  -4.98%  "
    "D49" = "3.77"
    "E49" = "  +0.76%  "
    "D50" = "21.08"
    "E50" = "  -0.15%  "
    "D51" = "0.624"
    "E51" = "  +1.28%  "
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $cell = $ws.Range($addr)

    $isPlainDecimal = $newValue -match '^-?\d+(\.\d+)?$'

    if ($isPlainDecimal) {
        # Force text storage, matching how this column is stored in the workbook,
        # then restore the default style so no stray number format sticks around.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    } else {
        $cell.Value = $newValue
    }
}
